# IEEEST supports remote bus specified by `busr`
#
# - Insert a new "busr" column into the IEEEST sheet, right after "MODE"
#   (i.e. before "A1"), and set its value for the existing data row.
# - Update the existing IEEEST row: fill in its uid, and bump MODE from 1 to 5
#   (to reflect remote-bus voltage-sensing mode).
# - Re-order the sheet tabs so IEEEST comes before Toggler.

$wb = $excel.ActiveWorkbook
$ieeest = $wb.Worksheets.Item("IEEEST")

# "MODE" is column F (6); insert the new "busr" column right after it,
# shifting A1..VCL one column to the right.
$ieeest.Columns.Item(7).Insert()

$ieeest.Cells.Item(1, 7).Value = "busr"
$ieeest.Cells.Item(2, 7).Value = 2

# Fill in the uid for the existing row, and update MODE to 5.
$ieeest.Cells.Item(2, 1).Value = 0
$ieeest.Cells.Item(2, 6).Value = 5

# Move the IEEEST tab so it sits before the Toggler tab.
$toggler = $wb.Worksheets.Item("Toggler")
$ieeest.Move($toggler)
